$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the missing PriceChange / UpDown values for the last existing
#     row (row 9) now that the next day's data is known ---
$ws.Range("X9").Value = -1.2099989999999963
$ws.Range("Y9").Value = "Down"

# --- Append the new trading-day row (row 10) ---
$ws.Range("A10").Value = 42653.87903935185
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = "Buy"
$ws.Range("D10").Value = 20
$ws.Range("E10").Value = 24940
$ws.Range("F10").Value = 2705
$ws.Range("G10").Value = 52
$ws.Range("H10").Value = 41
$ws.Range("I10").Value = 86
$ws.Range("J10").Value = 13
$ws.Range("K10").Value = 14682
$ws.Range("L10").Value = 350
$ws.Range("M10").Value = 276
$ws.Range("N10").Value = 109
$ws.Range("O10").Value = 17
$ws.Range("P10").Value = "Bag"
$ws.Range("Q10").Value = 35.550971360736582
$ws.Range("R10").Value = -24.44

# S10 / T10 use the percentage format already applied to S9 / T9 - copy the
# formatting across before writing the values so the existing style (s="2")
# is reused instead of a new one being minted.
$ws.Range("S9").Copy()
$ws.Range("S10").PasteSpecial(-4122)
$ws.Range("S10").Value = -0.1153

$ws.Range("T9").Copy()
$ws.Range("T10").PasteSpecial(-4122)
$ws.Range("T10").Value = -0.047

$ws.Range("U10").Value = 6.45
$ws.Range("V10").Value = 1.88
$ws.Range("W10").Value = 1

$excel.CutCopyMode = 0

# --- Nudge the "best fit" column widths slightly wider, matching the
#     auto-fit recalculation that happens when the new row is appended ---
$ws.Columns.Item(1).ColumnWidth = 14.541666666666666
$ws.Columns.Item(2).ColumnWidth = 7.666666666666667
$ws.Columns.Item(3).ColumnWidth = 5.666666666666667
$ws.Columns.Item(4).ColumnWidth = 11.291666666666666
$ws.Columns.Item(5).ColumnWidth = 8.666666666666666
$ws.Columns.Item(6).ColumnWidth = 11.416666666666666
$ws.Columns.Item(7).ColumnWidth = 18.416666666666668
$ws.Columns.Item(8).ColumnWidth = 18.541666666666668
$ws.Columns.Item(9).ColumnWidth = 19.541666666666668
$ws.Columns.Item(10).ColumnWidth = 19.791666666666668
$ws.Columns.Item(11).ColumnWidth = 9.541666666666666
$ws.Columns.Item(12).ColumnWidth = 13.541666666666666
$ws.Columns.Item(13).ColumnWidth = 13.791666666666666
